$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.769.09"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "3.702.95"
$ws.Range("E3").Value = "  +3.91%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.91"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +18.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "662.89"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("E8").Value = "  +4.74%  "

$ws.Range("E9").Value = "  +3.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").Value = "3.701.25"
$ws.Range("E11").Value = "  +3.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.58%  "

$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.57"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.04%  "

$ws.Range("D15").Value = "4.389.00"
$ws.Range("E15").Value = "  +3.85%  "

$ws.Range("D16").Value = "96.390.05"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000264"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +10.99%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.707.01"
$ws.Range("E19").Value = "  +3.32%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.547"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "514.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000212"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.47%  "

$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.04"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.46%  "

$ws.Range("E29").Value = "  +12.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("E33").Value = "  +1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.979"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.91%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.74"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.54%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "613.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "42.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +27.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.159"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.973"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +8.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.97%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.12"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0444"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.419"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +22.25%  "

$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.60"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.45%  "
